$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.302.85"
$ws.Range("E2").Value = "  +1.90%  "
$ws.Range("D3").Value = "1.648.94"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'217.76"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").Value = "'0.506"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").Value = "'19.98"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "1.878.57"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "1.673.01"
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "'63.57"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").Value = "26.297.07"
$ws.Range("E18").Value = "  +1.72%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "'196.76"
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("D22").Value = "'10.08"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("E24").Value = "  -2.09%  "
$ws.Range("D25").Value = "'143.17"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").Value = "'0.126"
$ws.Range("E27").Value = "  +1.65%  "
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("D29").Value = "'15.68"
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("D30").Value = "'1.26"
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("D31").Value = "'0.0505"
$ws.Range("E31").Value = "  +2.40%  "
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("E34").Value = "  +2.42%  "
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("D36").Value = "'0.917"
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("D37").Value = "'0.558"
$ws.Range("E37").Value = "  +2.27%  "
$ws.Range("D38").Value = "1.137.98"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").Value = "'2.49"
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").Value = "'5.65"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").Value = "'100.56"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("D45").Value = "1.788.14"
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("D46").Value = "'56.52"
$ws.Range("E46").Value = "  +2.05%  "
$ws.Range("D47").Value = "'1.50"
$ws.Range("E47").Value = "  +3.29%  "
$ws.Range("E48").Value = "  +3.14%  "
$ws.Range("D49").Value = "'7.73"
$ws.Range("E49").Value = "  +3.21%  "
$ws.Range("D50").Value = "'0.418"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("E51").Value = "  +2.25%  "
